$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear cells that are removed in the new layout ---
$ws.Range("B3:F3").ClearContents()
$ws.Range("C4:F4").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("C11:F11").ClearContents()
$ws.Range("A12:F12").ClearContents()
$ws.Range("B20").ClearContents()
$ws.Range("C20:F20").ClearContents()

# --- Update existing cells whose values changed ---
$ws.Range("C5").Value = "Spring 2022"
$ws.Range("C6").Value = "CPSC 4148"
$ws.Range("C7").Value = "CPSC 4155"
$ws.Range("A20").Value = "Fall 2023"

# --- Fall 2022 / Spring 2022 / Summer 2022 block (rows 3-15) ---
$ws.Range("E7").Value = "Summer 2022"

$ws.Range("C8").Value = "CPSC 4157"
$ws.Range("D8").Value = 3
$ws.Range("E8").Value = "CPSC 4176"
$ws.Range("F8").Value = 3

$ws.Range("C9").Value = "CPSC 4175"
$ws.Range("D9").Value = 3
$ws.Range("F9").Value = "Credits"

$ws.Range("C13").Value = "Total"
$ws.Range("D13").Formula = "=SUM(D6:D12)"

$ws.Range("E15").Value = "Total"
$ws.Range("F15").Formula = "=SUM(F8:F14)"

# --- Fall 2023 / Spring 2023 / Summer 2023 block (rows 20-32) ---
$ws.Range("B22").Value = "Credits"
$ws.Range("C22").Value = "Spring 2023"

$ws.Range("D24").Value = "Credits"
$ws.Range("E24").Value = "Summer 2023"

$ws.Range("F26").Value = "Credits"

$ws.Range("A28").Value = "Total"
$ws.Range("B28").Formula = "=SUM(B21:B27)"

$ws.Range("C30").Value = "Total"
$ws.Range("D30").Formula = "=SUM(D23:D29)"

$ws.Range("E32").Value = "Total"
$ws.Range("F32").Formula = "=SUM(F25:F31)"
